# The document's content/paragraphs are unchanged by this edit; the only
# structural change is where the hidden "_GoBack" bookmark (Word's
# "last edit location" marker) sits. It used to sit right after the
# "Purchasing Order" run (end of that list item). The new save moves it
# into the "Menjaga kebersihan dan kerapihan ..." paragraph, splitting
# the run "enjaga kebersihan dan kerapihan " into "enjaga kebersi" +
# bookmark + "han dan kerapihan " - i.e. the cursor was last sitting
# right after "...kebersi" when the document was saved.

$d = $word.ActiveDocument

# Locate the split point by searching for the text leading up to it,
# rather than hard-coding paragraph/character indices.
$target = $d.Content
$found = $target.Find.Execute("Menjaga kebersi", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate 'Menjaga kebersi' to reposition the _GoBack bookmark."
}

# Collapse the found range to its end point (right after "...kebersi").
$target.Collapse(0)

# Re-adding a bookmark with the same name moves it to the new range
# (Word bookmark names are unique within a document).
$d.Bookmarks.Add("_GoBack", $target)

Write-Output "Moved _GoBack bookmark to after 'Menjaga kebersi'."
